# Insert a new data row at row 294 (shifts existing rows 294-404 down to 295-405)
# and populate it with the new record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above current row 294; Excel shifts everything below
# down by one row and the worksheet dimension grows to A1:R405.
$ws.Rows.Item(294).Insert()

# Populate the newly inserted row 294 with the new record's values.
$ws.Cells.Item(294, 1).Value = 3
$ws.Cells.Item(294, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(294, 3).Value = 'Coquimbo'
$ws.Cells.Item(294, 4).Value = 44795
$ws.Cells.Item(294, 5).Value = 5
$ws.Cells.Item(294, 6).Value = 100112043
$ws.Cells.Item(294, 7).Value = 'Pepino ensalada'
$ws.Cells.Item(294, 8).Value = 'Sin especificar'
$ws.Cells.Item(294, 9).Value = 'Primera'
$ws.Cells.Item(294, 10).Value = 145
$ws.Cells.Item(294, 11).Value = 20000
$ws.Cells.Item(294, 12).Value = 21000
$ws.Cells.Item(294, 13).Value = 20517
$ws.Cells.Item(294, 14).Value = '$/caja 70 unidades'
$ws.Cells.Item(294, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(294, 16).Value = 293
$ws.Cells.Item(294, 17).Value = 70
$ws.Cells.Item(294, 18).Value = 'Hortaliza'
